# Slide 9 ("Grammar Rules Relevant to Arrays"), Content Placeholder 2:
# the "variable = ( varId | paramId) { indexExpr | fieldExpr } ." grammar
# rule is missing a space before the closing paren. Insert one so it reads
# "variable = ( varId | paramId ) { indexExpr | fieldExpr } .".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$paramIdPos = $fullText.IndexOf("paramId")

# 1-based character position immediately following "paramId" (i.e. the
# start of the existing ") { " run).
$afterParamId = $paramIdPos + 7 + 1

# The two characters ") " that currently open that run.
$closeParen = $tr.Characters($afterParamId, 2)

# Type a new " ) " right before them -- this is what splits a fresh run
# off of the ") { " run instead of merging into the "paramId" run.
$closeParen.InsertBefore(" ) ") | Out-Null

# The original ") " is now redundant (duplicated by the text we just
# inserted), so delete it, leaving behind "{ " in that run and our new
# " ) " run just after "paramId".
$fullText2 = $tr.Text
$paramIdPos2 = $fullText2.IndexOf("paramId")
$dupStart = $paramIdPos2 + 7 + 1 + 3
$dup = $tr.Characters($dupStart, 2)
$dup.Text = ""
